$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header cells right to make room for the new column:
# "Besar Gaji" moves from J1 -> K1, "Posisi" moves from I1 -> J1
# (Copy with a destination preserves both value and cell style.)
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("I1").Copy($ws.Range("J1"))

# Put the new header text in the freed-up I1 cell (keeps I1's existing style)
$ws.Range("I1").Value = "Jenis Instansi (Lokal/Nasional/Internasional)"

# Widen column I (9) to fit the new, longer header text
$ws.Columns.Item(9).ColumnWidth = 42.166666666666664

# Move the active cell/selection to K1, matching the new last header cell
$ws.Range("K1").Select()
